# Update "想去人数" (F) and "最低票价" (G) figures for the 展览 and 全部类型
# sheets, matching the refreshed data pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 6564
    $ws.Range("F3").Value = 188
    $ws.Range("F6").Value = 1981
    $ws.Range("G6").Value = 55
    $ws.Range("F7").Value = 1513
    $ws.Range("F8").Value = 307
    $ws.Range("F9").Value = 1005
    $ws.Range("F10").Value = 389
    $ws.Range("F12").Value = 5628
}
